$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# Row 42: I42 gets a commit-id value (new, wrapped text)
# ---------------------------------------------------------------
# (string created further below, in the exact order required to
#  reproduce the shared-string table ordering of the target file)

# ---------------------------------------------------------------
# Row 44: new data row (test 43)
# ---------------------------------------------------------------
$ws.Range("B44").Value() = 2
$ws.Range("C44").Value() = 100
$ws.Range("D44").Value() = 24
$ws.Range("E44").Value() = 1

# New shared strings must be created in this order so that the
# generated sharedStrings.xml matches the target ordering:
#   1) F44 text
#   2) I44 text
#   3) I42 text
#   4) H44 text
#   5) H45 text
$ws.Range("F44").Value() = "Базовая версия модели классификатора"
$ws.Range("F44").WrapText() = $true

$ws.Range("I44").Value() = "86b705a"

$ws.Range("I42").WrapText() = $true
$ws.Range("I42").Value() = "22e3876`n"

$ws.Range("H44").Value() = "Точность: Train: 83.5%, Val: 80.5%, Test: 63.3%. "

# G44 reuses an already existing shared string, keeps pre-existing style
$ws.Range("G44").Value() = "параметры теста 4"

$ws.Rows.Item(44).RowHeight() = 30

# ---------------------------------------------------------------
# Row 45: new "best result" data row (test 44), highlighted like row 40
# ---------------------------------------------------------------
$ws.Range("A40").Copy($ws.Range("A45"))
$ws.Range("A45").Value() = 44

$ws.Range("B40:E40").Copy($ws.Range("B45:E45"))
$ws.Range("B45").Value() = 2
$ws.Range("C45").Value() = 100
$ws.Range("D45").Value() = 24
$ws.Range("E45").Value() = 1

$ws.Range("F40").Copy($ws.Range("F45"))
$ws.Range("F45").Value() = "IoU/Dist = 20/80"

$ws.Range("G40:I40").Copy($ws.Range("G45:I45"))
$ws.Range("G45").Value() = "параметры теста 4"

$ws.Range("H45").Value() = "Точность: Train: 74.2%, Val: 63.7%, Test: 64%. "

$ws.Range("I45").ClearContents()

$ws.Rows.Item(45).RowHeight() = 30

# ---------------------------------------------------------------
# Update the active selection shown when the sheet is opened
# ---------------------------------------------------------------
$ws.Range("E49").Select()
